$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simatm40_1_1_3_JH")

# Add new test-case row value in A2 ("run")
$ws.Range("A2").Value = "run"

# Move the active selection from D17 to D2
$ws.Range("D2").Select()
